# Auto-generated edit script applying Atomos_Profits.xlsx diff
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ figures across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8264.556
$ws.Range("I19").Value = 297
$ws.Range("K19").Value = 297
$ws.Range("M19").Value = -122
$ws.Range("H28").Value = 634.25
$ws.Range("I28").Value = 146.17647
$ws.Range("J28").Value = 3400
$ws.Range("K28").Value = 146.17647
$ws.Range("L28").Value = 3400
$ws.Range("M28").Value = 338.82353
$ws.Range("N28").Value = -4370
$ws.Range("H76").Value = 4115
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4115
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 4115
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -4745
$ws.Range("H79").Value = 4115
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4115
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 4115
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -6299
$ws.Range("H113").Value = 5704.4287
$ws.Range("I113").Value = 4917.143
$ws.Range("J113").Value = 6491.7144
$ws.Range("K113").Value = 4917.143
$ws.Range("L113").Value = 6491.7144
$ws.Range("M113").Value = -1663.143
$ws.Range("N113").Value = -12999.7144
$ws.Range("H125").Value = 1857.1428
$ws.Range("I125").Value = 1166.6666
$ws.Range("J125").Value = 2375
$ws.Range("K125").Value = 10499.9994
$ws.Range("L125").Value = 21375
$ws.Range("M125").Value = -8039.999400000001
$ws.Range("N125").Value = -26295
$ws.Range("H137").Value = 2030.12
$ws.Range("I137").Value = 1346
$ws.Range("J137").Value = 2661.6155
$ws.Range("K137").Value = 4038
$ws.Range("L137").Value = 7984.8465
$ws.Range("M137").Value = -1488
$ws.Range("N137").Value = -13084.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6251872
$ws.Range("I2").Value = 15626016
$ws.Range("J2").Value = 2442.125
$ws.Range("K2").Value = 15626016
$ws.Range("L2").Value = 2442.125
$ws.Range("M2").Value = -15625903
$ws.Range("N2").Value = -2668.125
$ws.Range("H23").Value = 25667.5
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20518
$ws.Range("H32").Value = 3971.4524
$ws.Range("I32").Value = 2742.5615
$ws.Range("K32").Value = 2742.5615
$ws.Range("M32").Value = -2455.5615
$ws.Range("H45").Value = 1559.1389
$ws.Range("I45").Value = 991.09375
$ws.Range("K45").Value = 991.09375
$ws.Range("M45").Value = -614.09375
$ws.Range("H61").Value = 2774.2195
$ws.Range("I61").Value = 927.05554
$ws.Range("K61").Value = 927.05554
$ws.Range("M61").Value = -715.05554
$ws.Range("H63").Value = 4125.3125
$ws.Range("I63").Value = 2200.5
$ws.Range("K63").Value = 2200.5
$ws.Range("M63").Value = -1514.5
$ws.Range("H66").Value = 4125.3125
$ws.Range("I66").Value = 2200.5
$ws.Range("K66").Value = 11002.5
$ws.Range("M66").Value = -7570.5
$ws.Range("H97").Value = 900
$ws.Range("I97").Value = 875
$ws.Range("K97").Value = 875
$ws.Range("M97").Value = -379
$ws.Range("H116").Value = 6251872
$ws.Range("I116").Value = 15626016
$ws.Range("J116").Value = 2442.125
$ws.Range("K116").Value = 15626016
$ws.Range("L116").Value = 2442.125
$ws.Range("M116").Value = -15623722
$ws.Range("N116").Value = -7030.125
$ws.Range("H122").Value = 1760.25
$ws.Range("I122").Value = 1277.5834
$ws.Range("J122").Value = 2484.25
$ws.Range("K122").Value = 3832.7502
$ws.Range("L122").Value = 7452.75
$ws.Range("M122").Value = -1382.7502
$ws.Range("N122").Value = -12352.75
$ws.Range("H136").Value = 2774.2195
$ws.Range("I136").Value = 927.05554
$ws.Range("K136").Value = 2781.16662
$ws.Range("M136").Value = -231.16662

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6251872
$ws.Range("I3").Value = 15626016
$ws.Range("J3").Value = 2442.125
$ws.Range("K3").Value = 15626016
$ws.Range("L3").Value = 2442.125
$ws.Range("M3").Value = -15625902
$ws.Range("N3").Value = -2670.125
$ws.Range("H105").Value = 1637.0605
$ws.Range("I105").Value = 1318.9474
$ws.Range("J105").Value = 2068.7856
$ws.Range("K105").Value = 1318.9474
$ws.Range("L105").Value = 2068.7856
$ws.Range("M105").Value = 428.0526
$ws.Range("N105").Value = -5562.7856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2333.754
$ws.Range("I31").Value = 1439.7435
$ws.Range("J31").Value = 3674.7693
$ws.Range("K31").Value = 1439.7435
$ws.Range("L31").Value = 3674.7693
$ws.Range("M31").Value = -1144.7435
$ws.Range("N31").Value = -4264.7693
$ws.Range("H34").Value = 2333.754
$ws.Range("I34").Value = 1439.7435
$ws.Range("J34").Value = 3674.7693
$ws.Range("K34").Value = 1439.7435
$ws.Range("L34").Value = 3674.7693
$ws.Range("M34").Value = -1237.7435
$ws.Range("N34").Value = -4078.7693
$ws.Range("H58").Value = 20003218
$ws.Range("I58").Value = 1897.2858
$ws.Range("J58").Value = 45459444
$ws.Range("K58").Value = 1897.2858
$ws.Range("L58").Value = 45459444
$ws.Range("M58").Value = -1694.2858
$ws.Range("N58").Value = -45459850
$ws.Range("H94").Value = 3627.889
$ws.Range("I94").Value = 5691.75
$ws.Range("J94").Value = 1976.8
$ws.Range("K94").Value = 5691.75
$ws.Range("L94").Value = 1976.8
$ws.Range("M94").Value = -5240.75
$ws.Range("N94").Value = -2878.8
$ws.Range("H134").Value = 1826.9259
$ws.Range("I134").Value = 1418.5814
$ws.Range("J134").Value = 3423.182
$ws.Range("K134").Value = 4255.7442
$ws.Range("L134").Value = 10269.546
$ws.Range("M134").Value = -1720.7442
$ws.Range("N134").Value = -15339.546
$ws.Range("H136").Value = 20003218
$ws.Range("I136").Value = 1897.2858
$ws.Range("J136").Value = 45459444
$ws.Range("K136").Value = 5691.857400000001
$ws.Range("L136").Value = 136378332
$ws.Range("M136").Value = -3141.857400000001
$ws.Range("N136").Value = -136383432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 883.8732
$ws.Range("J107").Value = 1124.1945
$ws.Range("L107").Value = 3372.5835
$ws.Range("N107").Value = -7212.583500000001
$ws.Range("H121").Value = 42583.375
$ws.Range("I121").Value = 272.16666
$ws.Range("J121").Value = 169517
$ws.Range("K121").Value = 816.4999799999999
$ws.Range("L121").Value = 508551
$ws.Range("M121").Value = 493.5000200000001
$ws.Range("N121").Value = -511171
$ws.Range("H131").Value = 1737.6052
$ws.Range("I131").Value = 1578
$ws.Range("J131").Value = 1794.6072
$ws.Range("K131").Value = 4734
$ws.Range("L131").Value = 5383.821599999999
$ws.Range("M131").Value = 306
$ws.Range("N131").Value = -15463.8216
$ws.Range("H132").Value = 1422.2941
$ws.Range("I132").Value = 1195.2667
$ws.Range("K132").Value = 10757.4003
$ws.Range("M132").Value = -8227.400299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3932.0833
$ws.Range("I80").Value = 3486.4285
$ws.Range("J80").Value = 4556
$ws.Range("K80").Value = 3486.4285
$ws.Range("L80").Value = 4556
$ws.Range("M80").Value = -2488.4285
$ws.Range("N80").Value = -6552
$ws.Range("H83").Value = 3932.0833
$ws.Range("I83").Value = 3486.4285
$ws.Range("J83").Value = 4556
$ws.Range("K83").Value = 17432.1425
$ws.Range("L83").Value = 22780
$ws.Range("M83").Value = -12440.1425
$ws.Range("N83").Value = -32764
$ws.Range("H132").Value = 5024.9287
$ws.Range("I132").Value = 4024.8
$ws.Range("J132").Value = 5580.5557
$ws.Range("K132").Value = 12074.4
$ws.Range("L132").Value = 16741.6671
$ws.Range("M132").Value = -9544.400000000001
$ws.Range("N132").Value = -21801.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1635.9445
$ws.Range("I93").Value = 1307.909
$ws.Range("J93").Value = 2151.4285
$ws.Range("K93").Value = 1307.909
$ws.Range("L93").Value = 2151.4285
$ws.Range("M93").Value = -59.90900000000011
$ws.Range("N93").Value = -4647.4285
$ws.Range("H132").Value = 4135.636
$ws.Range("I132").Value = 2873.5
$ws.Range("K132").Value = 8620.5
$ws.Range("M132").Value = -6090.5
$ws.Range("H133").Value = 32260.8
$ws.Range("J133").Value = 32260.8
$ws.Range("L133").Value = 32260.8
$ws.Range("N133").Value = -37320.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 17667.834
$ws.Range("H33").Value = 6500
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 6500
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H126").Value = 4547387
$ws.Range("J126").Value = 7694616
$ws.Range("L126").Value = 23083848
$ws.Range("N126").Value = -23088788
